{"js": "/*\n * Replace the 100 arithmetic-expression strings in the single 20x5 table\n * (row-major order) with their updated values, per the commit diff.\n * Only the w:t text content of each cell is changed; run/paragraph\n * formatting (fonts, size, alignment) is left untouched because we set\n * TableCell.value instead of touching runs/fonts directly.\n */\nconst expectedOld = [\n  \"37+26=63\",\n  \"53-43=10\",\n  \"25+69=94\",\n  \"32+8=40\",\n  \"5+89=94\",\n  \"97-34=63\",\n  \"55-3=52\",\n  \"24-23=1\",\n  \"92-39=53\",\n  \"83-81=2\",\n  \"33-13=20\",\n  \"42+15=57\",\n  \"73-39=34\",\n  \"95+3=98\",\n  \"36-10=26\",\n  \"64-41=23\",\n  \"41+5=46\",\n  \"80-52=28\",\n  \"1+18=19\",\n  \"12-9=3\",\n  \"33+41=74\",\n  \"70-40=30\",\n  \"39+16=55\",\n  \"37+59=96\",\n  \"31-8=23\",\n  \"58+5=63\",\n  \"72-45=27\",\n  \"17+77=94\",\n  \"49-6=43\",\n  \"40-0=40\",\n  \"75+8=83\",\n  \"26+43=69\",\n  \"29+10=39\",\n  \"44-8=36\",\n  \"25+3=28\",\n  \"73-49=24\",\n  \"43+40=83\",\n  \"22-4=18\",\n  \"47+15=62\",\n  \"1+72=73\",\n  \"40-6=34\",\n  \"75-48=27\",\n  \"88+11=99\",\n  \"34+32=66\",\n  \"96-41=55\",\n  \"45+42=87\",\n  \"29+60=89\",\n  \"98-22=76\",\n  \"25+58=83\",\n  \"62+5=67\",\n  \"97-86=11\",\n  \"62+2=64\",\n  \"43+36=79\",\n  \"80-20=60\",\n  \"8+13=21\",\n  \"30-4=26\",\n  \"70-44=26\",\n  \"28+45=73\",\n  \"29-23=6\",\n  \"47-5=42\",\n  \"82-29=53\",\n  \"63-51=12\",\n  \"97-93=4\",\n  \"89-0=89\",\n  \"94-28=66\",\n  \"18-13=5\",\n  \"63+35=98\",\n  \"6+83=89\",\n  \"70-69=1\",\n  \"54-16=38\",\n  \"70+20=90\",\n  \"85-56=29\",\n  \"93-57=36\",\n  \"8-5=3\",\n  \"43+35=78\",\n  \"47+51=98\",\n  \"20+10=30\",\n  \"81-18=63\",\n  \"56+18=74\",\n  \"3+19=22\",\n  \"0+55=55\",\n  \"44+42=86\",\n  \"99-40=59\",\n  \"59+40=99\",\n  \"15+41=56\",\n  \"24-21=3\",\n  \"77+1=78\",\n  \"22+36=58\",\n  \"89-51=38\",\n  \"6+73=79\",\n  \"35+53=88\",\n  \"43-1=42\",\n  \"94-4=90\",\n  \"28+20=48\",\n  \"60-49=11\",\n  \"69-30=39\",\n  \"85-75=10\",\n  \"89-0=89\",\n  \"0+13=13\",\n  \"75+3=78\"\n];\nconst newValues = [\n  \"17+35=52\",\n  \"52+28=80\",\n  \"51+32=83\",\n  \"37+55=92\",\n  \"44+45=89\",\n  \"28-15=13\",\n  \"42+46=88\",\n  \"31+9=40\",\n  \"81-42=39\",\n  \"65-4=61\",\n  \"15+4=19\",\n  \"83+14=97\",\n  \"9+55=64\",\n  \"66-1=65\",\n  \"63-10=53\",\n  \"25+15=40\",\n  \"24+72=96\",\n  \"73-29=44\",\n  \"82-11=71\",\n  \"44-9=35\",\n  \"14+53=67\",\n  \"44+43=87\",\n  \"10+81=91\",\n  \"34+19=53\",\n  \"66-8=58\",\n  \"79-28=51\",\n  \"13-5=8\",\n  \"4+1=5\",\n  \"44-35=9\",\n  \"95-41=54\",\n  \"92-58=34\",\n  \"35+12=47\",\n  \"61-42=19\",\n  \"23+52=75\",\n  \"62+34=96\",\n  \"74-58=16\",\n  \"17+10=27\",\n  \"71-55=16\",\n  \"30-9=21\",\n  \"89-77=12\",\n  \"27-21=6\",\n  \"13+33=46\",\n  \"44-24=20\",\n  \"67-1=66\",\n  \"52+44=96\",\n  \"74+15=89\",\n  \"49-18=31\",\n  \"77-46=31\",\n  \"65-45=20\",\n  \"51+43=94\",\n  \"75-34=41\",\n  \"86-73=13\",\n  \"14+19=33\",\n  \"11+16=27\",\n  \"24+45=69\",\n  \"89-2=87\",\n  \"9+33=42\",\n  \"99-17=82\",\n  \"31+40=71\",\n  \"20+11=31\",\n  \"69-0=69\",\n  \"22+39=61\",\n  \"26+51=77\",\n  \"9+69=78\",\n  \"27+27=54\",\n  \"79-51=28\",\n  \"62-38=24\",\n  \"86+11=97\",\n  \"50-17=33\",\n  \"65-18=47\",\n  \"28+58=86\",\n  \"23+10=33\",\n  \"90-50=40\",\n  \"98-83=15\",\n  \"60-51=9\",\n  \"20+0=20\",\n  \"44+27=71\",\n  \"20+0=20\",\n  \"69-58=11\",\n  \"97-60=37\",\n  \"24+30=54\",\n  \"9+73=82\",\n  \"97-94=3\",\n  \"72+17=89\",\n  \"81-56=25\",\n  \"43+49=92\",\n  \"22-0=22\",\n  \"26+31=57\",\n  \"71+12=83\",\n  \"75-37=38\",\n  \"28-13=15\",\n  \"6+40=46\",\n  \"10+86=96\",\n  \"85+11=96\",\n  \"25+53=78\",\n  \"40+49=89\",\n  \"55-17=38\",\n  \"50+23=73\",\n  \"60+35=95\",\n  \"67-24=43\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No table found in document body\");\n}\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = table.values[0].length;\n\nif (rowCount * colCount !== newValues.length) {\n  throw new Error(\n    \"Table shape \" + rowCount + \"x\" + colCount +\n    \" does not match expected \" + newValues.length + \" cells\"\n  );\n}\n\n// Sanity-check (row-major) that the table still holds the pre-edit values\n// we expect before stomping on anything; two of the 100 original cells\n// happen to share the literal text \"89-0=89\" (rows 13 & 20), so matching is\n// done strictly by position, not by searching for old text.\nlet i = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const current = table.values[r][c];\n    if (current !== expectedOld[i]) {\n      throw new Error(\n        \"Cell (\" + r + \",\" + c + \") = \\\"\" + current +\n        \"\\\" does not match expected \\\"\" + expectedOld[i] + \"\\\" at position \" + i\n      );\n    }\n    i++;\n  }\n}\n\ni = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[i];\n    i++;\n  }\n}\n\nawait context.sync();\nreturn \"updated \" + i + \" cells\";\n", "ps1": "# Replace the 100 arithmetic-expression strings in the single 20x5 table\n# (row-major order) with their updated values, per the commit diff.\n# We set Cell.Range.Text directly so existing run formatting (fonts, size,\n# paragraph alignment) on each cell is preserved; only the text changes.\n\n$expectedOld = @(\"37+26=63\",\"53-43=10\",\"25+69=94\",\"32+8=40\",\"5+89=94\",\"97-34=63\",\"55-3=52\",\"24-23=1\",\"92-39=53\",\"83-81=2\",\"33-13=20\",\"42+15=57\",\"73-39=34\",\"95+3=98\",\"36-10=26\",\"64-41=23\",\"41+5=46\",\"80-52=28\",\"1+18=19\",\"12-9=3\",\"33+41=74\",\"70-40=30\",\"39+16=55\",\"37+59=96\",\"31-8=23\",\"58+5=63\",\"72-45=27\",\"17+77=94\",\"49-6=43\",\"40-0=40\",\"75+8=83\",\"26+43=69\",\"29+10=39\",\"44-8=36\",\"25+3=28\",\"73-49=24\",\"43+40=83\",\"22-4=18\",\"47+15=62\",\"1+72=73\",\"40-6=34\",\"75-48=27\",\"88+11=99\",\"34+32=66\",\"96-41=55\",\"45+42=87\",\"29+60=89\",\"98-22=76\",\"25+58=83\",\"62+5=67\",\"97-86=11\",\"62+2=64\",\"43+36=79\",\"80-20=60\",\"8+13=21\",\"30-4=26\",\"70-44=26\",\"28+45=73\",\"29-23=6\",\"47-5=42\",\"82-29=53\",\"63-51=12\",\"97-93=4\",\"89-0=89\",\"94-28=66\",\"18-13=5\",\"63+35=98\",\"6+83=89\",\"70-69=1\",\"54-16=38\",\"70+20=90\",\"85-56=29\",\"93-57=36\",\"8-5=3\",\"43+35=78\",\"47+51=98\",\"20+10=30\",\"81-18=63\",\"56+18=74\",\"3+19=22\",\"0+55=55\",\"44+42=86\",\"99-40=59\",\"59+40=99\",\"15+41=56\",\"24-21=3\",\"77+1=78\",\"22+36=58\",\"89-51=38\",\"6+73=79\",\"35+53=88\",\"43-1=42\",\"94-4=90\",\"28+20=48\",\"60-49=11\",\"69-30=39\",\"85-75=10\",\"89-0=89\",\"0+13=13\",\"75+3=78\")\n$newValues = @(\"17+35=52\",\"52+28=80\",\"51+32=83\",\"37+55=92\",\"44+45=89\",\"28-15=13\",\"42+46=88\",\"31+9=40\",\"81-42=39\",\"65-4=61\",\"15+4=19\",\"83+14=97\",\"9+55=64\",\"66-1=65\",\"63-10=53\",\"25+15=40\",\"24+72=96\",\"73-29=44\",\"82-11=71\",\"44-9=35\",\"14+53=67\",\"44+43=87\",\"10+81=91\",\"34+19=53\",\"66-8=58\",\"79-28=51\",\"13-5=8\",\"4+1=5\",\"44-35=9\",\"95-41=54\",\"92-58=34\",\"35+12=47\",\"61-42=19\",\"23+52=75\",\"62+34=96\",\"74-58=16\",\"17+10=27\",\"71-55=16\",\"30-9=21\",\"89-77=12\",\"27-21=6\",\"13+33=46\",\"44-24=20\",\"67-1=66\",\"52+44=96\",\"74+15=89\",\"49-18=31\",\"77-46=31\",\"65-45=20\",\"51+43=94\",\"75-34=41\",\"86-73=13\",\"14+19=33\",\"11+16=27\",\"24+45=69\",\"89-2=87\",\"9+33=42\",\"99-17=82\",\"31+40=71\",\"20+11=31\",\"69-0=69\",\"22+39=61\",\"26+51=77\",\"9+69=78\",\"27+27=54\",\"79-51=28\",\"62-38=24\",\"86+11=97\",\"50-17=33\",\"65-18=47\",\"28+58=86\",\"23+10=33\",\"90-50=40\",\"98-83=15\",\"60-51=9\",\"20+0=20\",\"44+27=71\",\"20+0=20\",\"69-58=11\",\"97-60=37\",\"24+30=54\",\"9+73=82\",\"97-94=3\",\"72+17=89\",\"81-56=25\",\"43+49=92\",\"22-0=22\",\"26+31=57\",\"71+12=83\",\"75-37=38\",\"28-13=15\",\"6+40=46\",\"10+86=96\",\"85+11=96\",\"25+53=78\",\"40+49=89\",\"55-17=38\",\"50+23=73\",\"60+35=95\",\"67-24=43\")\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\nif (($rows * $cols) -ne $newValues.Length) {\n    throw \"Table shape $rows x $cols does not match expected $($newValues.Length) cells\"\n}\n\n# Sanity-check (row-major) that the table still holds the pre-edit values we\n# expect before stomping on anything; two of the 100 original cells happen\n# to share the literal text \"89-0=89\" (rows 13 & 20), so matching is done\n# strictly by position, not by searching for old text.\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $current = $cell.Range.Text\n        $current = $current.TrimEnd([char]7).TrimEnd([char]13)\n        if ($current -ne $expectedOld[$i]) {\n            throw \"Cell ($r,$c) = `\"$current`\" does not match expected `\"$($expectedOld[$i])`\" at position $i\"\n        }\n        $i = $i + 1\n    }\n}\n\n$i = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n\nWrite-Output (\"updated \" + $i + \" cells\")\n"}
